$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6101797
$ws.Range("I40").Value = 2744.2144
$ws.Range("J40").Value = 9264269
$ws.Range("K40").Value = 2744.2144
$ws.Range("L40").Value = 9264269
$ws.Range("M40").Value = -2569.2144
$ws.Range("N40").Value = -9264619
$ws.Range("H58").Value = 372.2857
$ws.Range("I58").Value = 271.2
$ws.Range("J58").Value = 625
$ws.Range("K58").Value = 813.5999999999999
$ws.Range("L58").Value = 1875
$ws.Range("M58").Value = -663.5999999999999
$ws.Range("N58").Value = -2175
$ws.Range("H98").Value = 874.2
$ws.Range("I98").Value = 842.75
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 842.75
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 655.25
$ws.Range("N98").Value = -3996
$ws.Range("H112").Value = 2569.8572
$ws.Range("I112").Value = 1896.2
$ws.Range("K112").Value = 5688.6
$ws.Range("M112").Value = -4580.6
$ws.Range("H122").Value = 874.2
$ws.Range("I122").Value = 842.75
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2528.25
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -78.25
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 3051.8572
$ws.Range("I132").Value = 3051.8572
$ws.Range("K132").Value = 9155.571599999999
$ws.Range("M132").Value = -6625.571599999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7683.4185
$ws.Range("I32").Value = 7628.2617
$ws.Range("K32").Value = 7628.2617
$ws.Range("M32").Value = -7341.2617
$ws.Range("H45").Value = 2238.074
$ws.Range("I45").Value = 1715.9412
$ws.Range("K45").Value = 1715.9412
$ws.Range("M45").Value = -1338.9412
$ws.Range("H46").Value = 7984.1
$ws.Range("I46").Value = 5047.6665
$ws.Range("J46").Value = 9242.571
$ws.Range("K46").Value = 5047.6665
$ws.Range("L46").Value = 9242.571
$ws.Range("M46").Value = -4728.6665
$ws.Range("N46").Value = -9880.571
$ws.Range("H61").Value = 3668.7778
$ws.Range("I61").Value = 3590.4707
$ws.Range("K61").Value = 3590.4707
$ws.Range("M61").Value = -3378.4707
$ws.Range("H74").Value = 1450.9395
$ws.Range("I74").Value = 1093.2903
$ws.Range("K74").Value = 1093.2903
$ws.Range("M74").Value = -219.2902999999999
$ws.Range("H77").Value = 1450.9395
$ws.Range("I77").Value = 1093.2903
$ws.Range("K77").Value = 5466.451499999999
$ws.Range("M77").Value = -1098.451499999999
$ws.Range("H88").Value = 2362.1428
$ws.Range("J88").Value = 2563.4285
$ws.Range("L88").Value = 2563.4285
$ws.Range("N88").Value = -3375.4285
$ws.Range("H91").Value = 2362.1428
$ws.Range("J91").Value = 2563.4285
$ws.Range("L91").Value = 2563.4285
$ws.Range("N91").Value = -5371.4285
$ws.Range("H110").Value = 4631.1665
$ws.Range("I110").Value = 510.57144
$ws.Range("K110").Value = 510.57144
$ws.Range("M110").Value = 1534.42856
$ws.Range("H122").Value = 2911
$ws.Range("I122").Value = 2911
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8733
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6283
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2795.889
$ws.Range("J132").Value = 2879.25
$ws.Range("L132").Value = 8637.75
$ws.Range("N132").Value = -13697.75
$ws.Range("H136").Value = 3668.7778
$ws.Range("I136").Value = 3590.4707
$ws.Range("K136").Value = 10771.4121
$ws.Range("M136").Value = -8221.4121

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7197.706
$ws.Range("I20").Value = 6638.5835
$ws.Range("K20").Value = 6638.5835
$ws.Range("M20").Value = -6391.5835
$ws.Range("H22").Value = 999.6667
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H94").Value = 674.8333
$ws.Range("I94").Value = 671.9375
$ws.Range("J94").Value = 698
$ws.Range("K94").Value = 671.9375
$ws.Range("L94").Value = 698
$ws.Range("M94").Value = -220.9375
$ws.Range("N94").Value = -1600
$ws.Range("H102").Value = 11943.5
$ws.Range("I102").Value = 11943.5
$ws.Range("K102").Value = 11943.5
$ws.Range("M102").Value = -8698.5
$ws.Range("H134").Value = 8183.6
$ws.Range("I134").Value = 7731.25
$ws.Range("K134").Value = 23193.75
$ws.Range("M134").Value = -20658.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3130.8667
$ws.Range("I62").Value = 3107.889
$ws.Range("K62").Value = 3107.889
$ws.Range("M62").Value = -2483.889
$ws.Range("H65").Value = 3130.8667
$ws.Range("I65").Value = 3107.889
$ws.Range("K65").Value = 15539.445
$ws.Range("M65").Value = -12419.445
$ws.Range("H99").Value = 5879.381
$ws.Range("I99").Value = 6272.8335
$ws.Range("K99").Value = 6272.8335
$ws.Range("M99").Value = -4774.8335
$ws.Range("H107").Value = 1379.8334
$ws.Range("I107").Value = 509.0909
$ws.Range("K107").Value = 509.0909
$ws.Range("M107").Value = 1410.9091
$ws.Range("H122").Value = 2475.516
$ws.Range("I122").Value = 2664.92
$ws.Range("K122").Value = 7994.76
$ws.Range("M122").Value = -5544.76
$ws.Range("H126").Value = 5879.381
$ws.Range("I126").Value = 6272.8335
$ws.Range("K126").Value = 18818.5005
$ws.Range("M126").Value = -16348.5005
$ws.Range("H132").Value = 7537.5
$ws.Range("I132").Value = 5050
$ws.Range("K132").Value = 15150
$ws.Range("M132").Value = -12620
$ws.Range("H134").Value = 7275
$ws.Range("I134").Value = 4140
$ws.Range("J134").Value = 12500
$ws.Range("K134").Value = 12420
$ws.Range("L134").Value = 37500
$ws.Range("M134").Value = -9885
$ws.Range("N134").Value = -42570

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 28560292
$ws.Range("I4").Value = 26730380
$ws.Range("K4").Value = 80191140
$ws.Range("M4").Value = -80191028
$ws.Range("H7").Value = 484.57144
$ws.Range("I7").Value = 515.1667
$ws.Range("K7").Value = 1545.5001
$ws.Range("M7").Value = -1433.5001
$ws.Range("H12").Value = 857
$ws.Range("J12").Value = 979.2857
$ws.Range("L12").Value = 2937.8571
$ws.Range("N12").Value = -3283.8571
$ws.Range("H107").Value = 435473.34
$ws.Range("I107").Value = 458.45456
$ws.Range("J107").Value = 834237
$ws.Range("K107").Value = 1375.36368
$ws.Range("L107").Value = 2502711
$ws.Range("M107").Value = 544.6363200000001
$ws.Range("N107").Value = -2506551
$ws.Range("H129").Value = 555098.8
$ws.Range("J129").Value = 673353.4
$ws.Range("L129").Value = 2020060.2
$ws.Range("N129").Value = -2030060.2
$ws.Range("H132").Value = 1573.8158
$ws.Range("J132").Value = 1573.8158
$ws.Range("L132").Value = 14164.3422
$ws.Range("N132").Value = -19224.3422

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2692.639
$ws.Range("I102").Value = 1895.7812
$ws.Range("J102").Value = 9067.5
$ws.Range("K102").Value = 1895.7812
$ws.Range("L102").Value = 9067.5
$ws.Range("M102").Value = -273.7811999999999
$ws.Range("N102").Value = -12311.5
$ws.Range("H113").Value = 365022.28
$ws.Range("I113").Value = 445824.44
$ws.Range("K113").Value = 445824.44
$ws.Range("M113").Value = -443654.44
$ws.Range("H122").Value = 1756.95
$ws.Range("I122").Value = 1423.3077
$ws.Range("K122").Value = 4269.9231
$ws.Range("M122").Value = -1819.9231
$ws.Range("H126").Value = 5776.8887
$ws.Range("I126").Value = 4832
$ws.Range("K126").Value = 14496
$ws.Range("M126").Value = -12026
$ws.Range("H132").Value = 5495.5
$ws.Range("I132").Value = 2500
$ws.Range("K132").Value = 7500
$ws.Range("M132").Value = -4970

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1064.0667
$ws.Range("J82").Value = 1089.8
$ws.Range("L82").Value = 1089.8
$ws.Range("N82").Value = -1811.8
$ws.Range("H85").Value = 1064.0667
$ws.Range("J85").Value = 1089.8
$ws.Range("L85").Value = 1089.8
$ws.Range("N85").Value = -3585.8
$ws.Range("H100").Value = 6264.8823
$ws.Range("I100").Value = 5150.9
$ws.Range("J100").Value = 7856.2856
$ws.Range("K100").Value = 5150.9
$ws.Range("L100").Value = 7856.2856
$ws.Range("M100").Value = -4609.9
$ws.Range("N100").Value = -8938.285599999999
$ws.Range("H122").Value = 2111.8572
$ws.Range("I122").Value = 2111.8572
$ws.Range("K122").Value = 6335.571599999999
$ws.Range("M122").Value = -3885.571599999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2897
$ws.Range("I96").Value = 2947.5
$ws.Range("K96").Value = 2947.5
$ws.Range("M96").Value = -1574.5
$ws.Range("H100").Value = 927.5484
$ws.Range("I100").Value = 836.5789
$ws.Range("J100").Value = 1071.5834
$ws.Range("K100").Value = 1673.1578
$ws.Range("L100").Value = 2143.1668
$ws.Range("M100").Value = -1132.1578
$ws.Range("N100").Value = -3225.1668
$ws.Range("H126").Value = 4441.8203
$ws.Range("I126").Value = 4269.1714
$ws.Range("J126").Value = 5952.5
$ws.Range("K126").Value = 12807.5142
$ws.Range("L126").Value = 17857.5
$ws.Range("M126").Value = -10337.5142
$ws.Range("N126").Value = -22797.5
$ws.Range("H137").Value = 59999.5
$ws.Range("J137").Value = 59999.5
$ws.Range("L137").Value = 59999.5
$ws.Range("N137").Value = -70199.5
